$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
# Row 11
$ws.Range("H11").Value = 24.846153
$ws.Range("I11").Value = 24.846153
$ws.Range("K11").Value = 24.846153
$ws.Range("M11").Value = 115.153847

# Row 19
$ws.Range("H19").Value = 2163.4
$ws.Range("I19").Value = 0
$ws.Range("J19").Value = 2163.4
$ws.Range("K19").Value = 0
$ws.Range("L19").Value = 2163.4
$ws.Range("M19").ClearContents()
$ws.Range("N19").Value = -2513.4

# Row 33
$ws.Range("H33").Value = 1097.5883
$ws.Range("I33").Value = 686.8333
$ws.Range("K33").Value = 686.8333
$ws.Range("M33").Value = -457.8333

# Row 43
$ws.Range("H43").Value = 2740.6
$ws.Range("I43").Value = 2520.238
$ws.Range("J43").Value = 3254.7778
$ws.Range("K43").Value = 2520.238
$ws.Range("L43").Value = 3254.7778
$ws.Range("M43").Value = -2451.238
$ws.Range("N43").Value = -3392.7778

# Row 86
$ws.Range("H86").Value = 4285.3184
$ws.Range("I86").Value = 3271.9092
$ws.Range("K86").Value = 3271.9092
$ws.Range("M86").Value = -2148.9092

# Row 89
$ws.Range("H89").Value = 4285.3184
$ws.Range("I89").Value = 3271.9092
$ws.Range("K89").Value = 16359.546
$ws.Range("M89").Value = -10743.546

# Row 94
$ws.Range("H94").Value = 628
$ws.Range("I94").Value = 646.4286
$ws.Range("J94").Value = 499
$ws.Range("K94").Value = 646.4286
$ws.Range("L94").Value = 499
$ws.Range("M94").Value = -195.4286
$ws.Range("N94").Value = -1401

# Row 127
$ws.Range("H127").Value = 2052.5293
$ws.Range("I127").Value = 764.9231
$ws.Range("K127").Value = 2294.7693
$ws.Range("M127").Value = 2665.2307

# Row 137
$ws.Range("H137").Value = 3512.6978
$ws.Range("J137").Value = 4364.154
$ws.Range("L137").Value = 13092.462
$ws.Range("N137").Value = -18192.462

$ws = $wb.Worksheets.Item("ARM")
# Row 61
$ws.Range("H61").Value = 8976.895
$ws.Range("I61").Value = 8586.723
$ws.Range("K61").Value = 8586.723
$ws.Range("M61").Value = -8374.723

# Row 63
$ws.Range("H63").Value = 5599.8
$ws.Range("I63").Value = 5999.75
$ws.Range("K63").Value = 5999.75
$ws.Range("M63").Value = -5313.75

# Row 66
$ws.Range("H66").Value = 5599.8
$ws.Range("I66").Value = 5999.75
$ws.Range("K66").Value = 29998.75
$ws.Range("M66").Value = -26566.75

# Row 80
$ws.Range("H80").Value = 0
$ws.Range("J80").Value = 0
$ws.Range("L80").Value = 0
$ws.Range("N80").ClearContents()

# Row 83
$ws.Range("H83").Value = 0
$ws.Range("J83").Value = 0
$ws.Range("L83").Value = 0
$ws.Range("N83").ClearContents()

# Row 126
$ws.Range("H126").Value = 7988.3335
$ws.Range("I126").Value = 7988.3335
$ws.Range("K126").Value = 23965.0005
$ws.Range("M126").Value = -21495.0005

# Row 136
$ws.Range("H136").Value = 8976.895
$ws.Range("I136").Value = 8586.723
$ws.Range("K136").Value = 25760.169
$ws.Range("M136").Value = -23210.169

$ws = $wb.Worksheets.Item("BSM")
# Row 22
$ws.Range("H22").Value = 805.92
$ws.Range("I22").Value = 632.45
$ws.Range("K22").Value = 632.45
$ws.Range("M22").Value = -459.45

# Row 75
$ws.Range("H75").Value = 11957
$ws.Range("I75").Value = 12616.5
$ws.Range("J75").Value = 8000
$ws.Range("K75").Value = 12616.5
$ws.Range("L75").Value = 8000
$ws.Range("M75").Value = -11680.5
$ws.Range("N75").Value = -9872

# Row 78
$ws.Range("H78").Value = 11957
$ws.Range("I78").Value = 12616.5
$ws.Range("J78").Value = 8000
$ws.Range("K78").Value = 37849.5
$ws.Range("L78").Value = 24000
$ws.Range("M78").Value = -33169.5
$ws.Range("N78").Value = -33360

# Row 99
$ws.Range("H99").Value = 3777.8333
$ws.Range("I99").Value = 1844
$ws.Range("J99").Value = 4744.75
$ws.Range("K99").Value = 1844
$ws.Range("L99").Value = 4744.75
$ws.Range("M99").Value = -346
$ws.Range("N99").Value = -7740.75

# Row 105
$ws.Range("H105").Value = 7532.857
$ws.Range("I105").Value = 7788.5
$ws.Range("J105").Value = 5999
$ws.Range("K105").Value = 7788.5
$ws.Range("L105").Value = 5999
$ws.Range("M105").Value = -6041.5
$ws.Range("N105").Value = -9493

# Row 134
$ws.Range("H134").Value = 4892.849
$ws.Range("I134").Value = 4118
$ws.Range("J134").Value = 7851.364
$ws.Range("K134").Value = 12354
$ws.Range("L134").Value = 23554.092
$ws.Range("M134").Value = -9819
$ws.Range("N134").Value = -28624.092

$ws = $wb.Worksheets.Item("CRP")
# Row 39
$ws.Range("H39").Value = 9810.75
$ws.Range("I39").Value = 5498.143
$ws.Range("K39").Value = 5498.143
$ws.Range("M39").Value = -5107.143

# Row 47
$ws.Range("H47").Value = 26267
$ws.Range("J47").Value = 30022.666
$ws.Range("L47").Value = 30022.666
$ws.Range("N47").Value = -31154.666

# Row 48
$ws.Range("H48").Value = 51750
$ws.Range("J48").Value = 51750
$ws.Range("L48").Value = 51750
$ws.Range("N48").Value = -52702

# Row 49
$ws.Range("H49").Value = 9810.75
$ws.Range("I49").Value = 5498.143
$ws.Range("K49").Value = 5498.143
$ws.Range("M49").Value = -5316.143

# Row 86
$ws.Range("H86").Value = 5044.75
$ws.Range("I86").Value = 4432.75
$ws.Range("K86").Value = 4432.75
$ws.Range("M86").Value = -3309.75

# Row 89
$ws.Range("H89").Value = 5044.75
$ws.Range("I89").Value = 4432.75
$ws.Range("K89").Value = 22163.75
$ws.Range("M89").Value = -16547.75

# Row 99
$ws.Range("H99").Value = 24479.6
$ws.Range("I99").Value = 28424.75
$ws.Range("K99").Value = 28424.75
$ws.Range("M99").Value = -26926.75

# Row 102
$ws.Range("H102").Value = 94620.5
$ws.Range("J102").Value = 94620.5
$ws.Range("L102").Value = 94620.5
$ws.Range("N102").Value = -99488.5

# Row 126
$ws.Range("H126").Value = 24479.6
$ws.Range("I126").Value = 28424.75
$ws.Range("K126").Value = 85274.25
$ws.Range("M126").Value = -82804.25

# Row 141
$ws.Range("H141").Value = 250145.72
$ws.Range("J141").Value = 282688.53
$ws.Range("L141").Value = 282688.53
$ws.Range("N141").Value = -293048.53

$ws = $wb.Worksheets.Item("CUL")
# Row 12
$ws.Range("H12").Value = 126
$ws.Range("I12").Value = 10.5
$ws.Range("J12").Value = 149.1
$ws.Range("K12").Value = 31.5
$ws.Range("L12").Value = 447.3
$ws.Range("M12").Value = 141.5
$ws.Range("N12").Value = -793.3

# Row 115
$ws.Range("H115").Value = 12799.75
$ws.Range("I115").Value = 399.66666
$ws.Range("K115").Value = 1198.99998
$ws.Range("M115").Value = -23.99998000000005

# Row 122
$ws.Range("H122").Value = 2581.1667
$ws.Range("I122").Value = 5243.5
$ws.Range("J122").Value = 1250
$ws.Range("K122").Value = 47191.5
$ws.Range("L122").Value = 11250
$ws.Range("M122").Value = -44741.5
$ws.Range("N122").Value = -16150

$ws = $wb.Worksheets.Item("GSM")
# Row 70
$ws.Range("H70").Value = 5162.6
$ws.Range("I70").Value = 4422.6665
$ws.Range("J70").Value = 5293.1763
$ws.Range("K70").Value = 4422.6665
$ws.Range("L70").Value = 5293.1763
$ws.Range("M70").Value = -4152.6665
$ws.Range("N70").Value = -5833.1763

# Row 73
$ws.Range("H73").Value = 5162.6
$ws.Range("I73").Value = 4422.6665
$ws.Range("J73").Value = 5293.1763
$ws.Range("K73").Value = 4422.6665
$ws.Range("L73").Value = 5293.1763
$ws.Range("M73").Value = -3486.6665
$ws.Range("N73").Value = -7165.1763

# Row 97
$ws.Range("H97").Value = 982.75
$ws.Range("I97").Value = 974.125
$ws.Range("J97").Value = 1000
$ws.Range("K97").Value = 974.125
$ws.Range("L97").Value = 1000
$ws.Range("M97").Value = -478.125
$ws.Range("N97").Value = -1992

# Row 122
$ws.Range("H122").Value = 6923
$ws.Range("I122").Value = 4176.875
$ws.Range("K122").Value = 12530.625
$ws.Range("M122").Value = -10080.625

$ws = $wb.Worksheets.Item("LTW")
# Row 16
$ws.Range("H16").Value = 628.2759
$ws.Range("I16").Value = 536.6087
$ws.Range("J16").Value = 979.6667
$ws.Range("K16").Value = 536.6087
$ws.Range("L16").Value = 979.6667
$ws.Range("M16").Value = -366.6087
$ws.Range("N16").Value = -1319.6667

# Row 68
$ws.Range("H68").Value = 6808.643
$ws.Range("I68").Value = 5121.909
$ws.Range("K68").Value = 5121.909
$ws.Range("M68").Value = -4372.909

# Row 71
$ws.Range("H71").Value = 6808.643
$ws.Range("I71").Value = 5121.909
$ws.Range("K71").Value = 25609.545
$ws.Range("M71").Value = -21865.545

# Row 100
$ws.Range("H100").Value = 5499.5415
$ws.Range("I100").Value = 5226.591
$ws.Range("K100").Value = 5226.591
$ws.Range("M100").Value = -4685.591

$ws = $wb.Worksheets.Item("WVR")
# Row 11
$ws.Range("H11").Value = 22502.5
$ws.Range("I11").Value = 0
$ws.Range("K11").Value = 0
$ws.Range("M11").ClearContents()

# Row 37
$ws.Range("H37").Value = 62423.625
$ws.Range("J37").Value = 77078.39999999999
$ws.Range("L37").Value = 77078.39999999999
$ws.Range("N37").Value = -77484.39999999999

# Row 75
$ws.Range("H75").Value = 91666.664
$ws.Range("I75").Value = 95000
$ws.Range("J75").Value = 90000
$ws.Range("K75").Value = 95000
$ws.Range("L75").Value = 90000
$ws.Range("M75").Value = -94064
$ws.Range("N75").Value = -91872

# Row 78
$ws.Range("H78").Value = 91666.664
$ws.Range("I78").Value = 95000
$ws.Range("J78").Value = 90000
$ws.Range("K78").Value = 285000
$ws.Range("L78").Value = 270000
$ws.Range("M78").Value = -280320
$ws.Range("N78").Value = -279360

# Row 107
$ws.Range("H107").Value = 3667.1052
$ws.Range("I107").Value = 2962
$ws.Range("J107").Value = 6311.25
$ws.Range("K107").Value = 8886
$ws.Range("L107").Value = 18933.75
$ws.Range("M107").Value = -6966
$ws.Range("N107").Value = -22773.75

# Row 136
$ws.Range("H136").Value = 3668.6875
$ws.Range("I136").Value = 3338.077
$ws.Range("K136").Value = 10014.231
$ws.Range("M136").Value = -7464.231
